$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Louise, Joséphine - Arrivée
$ws.Range("A8").Value = 43343.95301059028
$ws.Range("C8").Value = "Louise, Joséphine"
$ws.Range("D8").Value = "Arrivée"
$ws.Range("E8").Value = 0.32291666666424135

# Row 9: Joséphine - Départ
$ws.Range("A9").Value = 43343.953370034724
$ws.Range("C9").Value = "Joséphine"
$ws.Range("D9").Value = "Départ"
$ws.Range("E9").Value = 0.7083333333357587
$ws.Range("F9").Value = "Goûter"
$ws.Range("G9").Value = "AR école"

# Row 10: Joséphine - Arrivée
$ws.Range("A10").Value = 43346.79166195602
$ws.Range("C10").Value = "Joséphine"
$ws.Range("D10").Value = "Arrivée"
$ws.Range("E10").Value = 0.38541666666424135

# Row 11: Joséphine - Départ
$ws.Range("A11").Value = 43346.79189278935
$ws.Range("C11").Value = "Joséphine"
$ws.Range("D11").Value = "Départ"
$ws.Range("E11").Value = 0.7083333333357587

# Apply formats matching existing rows via copy/paste-special (reuses existing style indices)
$ws.Range("A2").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)
$ws.Range("C2:D2").Copy()
$ws.Range("C8:D8").PasteSpecial(-4122)
$ws.Range("C9:D9").PasteSpecial(-4122)
$ws.Range("C10:D10").PasteSpecial(-4122)
$ws.Range("C11:D11").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E8:E11").PasteSpecial(-4122)
$ws.Range("F5:G5").Copy()
$ws.Range("F9:G9").PasteSpecial(-4122)
